$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Highlight existing "Connection/Variable Name" header-adjacent rows in yellow ---
# (rows 2,3 under the first table; rows 15-18 under the digital-pin table)
$ws.Range("A2:C2").Interior.Color = 65535
$ws.Range("A3:C3").Interior.Color = 65535
$ws.Range("A15:C15").Interior.Color = 65535
$ws.Range("A16:C16").Interior.Color = 65535
$ws.Range("A17:C17").Interior.Color = 65535
$ws.Range("A18:C18").Interior.Color = 65535

# --- New "partial reflection" multiplexer wiring block (rows 59-66) ---
$ws.Range("I59").Value = "vss"
$ws.Range("J59").Value = "x"
$ws.Range("K59").Value = "x"
$ws.Range("L59").Value = "c"

$ws.Range("I60").Value = "vee"
$ws.Range("J60").Value = "x"
$ws.Range("K60").Value = "x"
$ws.Range("L60").Value = "b"

$ws.Range("I61").Value = "inh"
$ws.Range("J61").Value = "x"
$ws.Range("K61").Value = "x"
$ws.Range("L61").Value = "a"

$ws.Range("I62").Value = "s5"
$ws.Range("J62").Value = "x"
$ws.Range("K62").Value = "x"
$ws.Range("L62").Value = "s3"

$ws.Range("F63").Value = "s7"
$ws.Range("G63").Value = "g"
$ws.Range("H63").Value = 5
$ws.Range("I63").Value = "s7"
$ws.Range("J63").Value = "x"
$ws.Range("K63").Value = "x"
$ws.Range("L63").Value = "s0"
$ws.Range("M63").Value = "s3"
$ws.Range("N63").Value = "g"
$ws.Range("O63").Value = 5

$ws.Range("F64").Value = "s6"
$ws.Range("G64").Value = "g"
$ws.Range("H64").Value = 5
$ws.Range("I64").Value = "com"
$ws.Range("J64").Value = "x"
$ws.Range("K64").Value = "x"
$ws.Range("L64").Value = "s1"
$ws.Range("M64").Value = "s2"
$ws.Range("N64").Value = "g"
$ws.Range("O64").Value = 5

$ws.Range("F65").Value = "s5"
$ws.Range("G65").Value = "g"
$ws.Range("H65").Value = 5
$ws.Range("I65").Value = "s6"
$ws.Range("J65").Value = "x"
$ws.Range("K65").Value = "x"
$ws.Range("L65").Value = "s2"
$ws.Range("M65").Value = "s1"
$ws.Range("N65").Value = "g"
$ws.Range("O65").Value = 5

$ws.Range("F66").Value = "s4"
$ws.Range("G66").Value = "g"
$ws.Range("H66").Value = 5
$ws.Range("I66").Value = "s4"
$ws.Range("J66").Value = "x"
$ws.Range("K66").Value = "x"
$ws.Range("L66").Value = "vdd"
$ws.Range("M66").Value = "s0"
$ws.Range("N66").Value = "g"
$ws.Range("O66").Value = 5

# --- Restore selection to match the saved view ---
[void]$ws.Range("E12").Select()
